# The post "「このイノセントな目に騙されるな…」" (row 310) was removed from the
# spreadsheet. Delete its entire row; Excel will automatically shift all
# subsequent rows up by one and adjust the used range/dimension.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(310).Delete()
